$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.458.23"
$ws.Range("E2").Value = "  -3.15%  "
$ws.Range("D3").Value = "2.469.87"
$ws.Range("E3").Value = "  -2.34%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.92"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.17%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.551"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.68%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  -4.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.65"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.97%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0782"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.01%  "
$ws.Range("E12").Value = "  -0.89%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.04"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.18%  "
$ws.Range("D14").Value = "2.850.54"
$ws.Range("E14").Value = "  -2.27%  "
$ws.Range("D15").Value = "2.456.68"
$ws.Range("E15").Value = "  -3.52%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.88"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.96%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.788"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.40%  "
$ws.Range("D18").Value = "41.403.12"
$ws.Range("E18").Value = "  -3.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.35"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.83%  "
$ws.Range("D20").Value = "0.0₃0924"
$ws.Range("E20").Value = "  -3.10%  "
$ws.Range("E21").Value = "  -8.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.64"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.81%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.71"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.41%  "
$ws.Range("E24").Value = "  -4.01%  "
$ws.Range("B25").Value = "ImmutableX"
$ws.Range("C25").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.91"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -6.17%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.25"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.91%  "
$ws.Range("E28").Value = "  -5.49%  "
$ws.Range("E29").Value = "  -5.37%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.84"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.86%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "151.98"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.00%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.51"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.23%  "
$ws.Range("E33").Value = "  -4.97%  "
$ws.Range("E34").Value = "  -2.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0747"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.46%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.07"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.74%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.90"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.22%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "17.08"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.81%  "
$ws.Range("E39").Value = "  -2.69%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.30"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.17%  "
$ws.Range("E41").Value = "  -7.57%  "
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.04"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.90%  "
$ws.Range("D44").Value = "1.991.61"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0287"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.19%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.06"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.80"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.41%  "
$ws.Range("D48").Value = "2.709.25"
$ws.Range("E48").Value = "  -2.21%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "70.06"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.40%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "97.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.95%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "75.13"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.56%  "
